$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A3").Value = 112442490
$ws.Range("B3").Value = 96720
$ws.Range("E3").Value = 220787
$ws.Range("Q3").Value = 518340
$ws.Range("R3").Value = 6608985
$ws.Range("S3").Value = 10

# Text cells
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("P3").Value = "Ålkilen, Vstm"
$ws.Range("T3").Value = "Örebro"
$ws.Range("U3").Value = "Lindesberg"
$ws.Range("V3").Value = "Västmanland"
$ws.Range("W3").Value = "Linde"
$ws.Range("AW3").Value = "Alexander Singer"
$ws.Range("AX3").Value = "Alexander Singer"

# Date-like text cells that must remain plain text (not auto-converted to Excel dates)
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-09-01"
$ws.Range("Y3").Style = "Normal"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-09-30"
$ws.Range("AA3").Style = "Normal"

# Boolean cells
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# Blank placeholder cells (present in the row but empty, matching source data)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Style = "Normal"

$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Style = "Normal"

$ws.Range("AY3").NumberFormat = "@"
$ws.Range("AY3").Style = "Normal"
